$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.328.51"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "3.102.76"
$ws.Range("E3").Value = "  -3.54%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.47"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.77"
$ws.Range("E6").Value = "  +3.58%  "
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("D8").Value = "3.102.64"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.91"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("E12").Value = "  +0.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.65"
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000242"
$ws.Range("E14").Value = "  -2.00%  "
$ws.Range("D15").Value = "3.616.94"
$ws.Range("E15").Value = "  -3.42%  "
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.21"
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("D18").Value = "63.865.67"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").Value = "3.102.81"
$ws.Range("E19").Value = "  -2.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "468.14"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.80"
$ws.Range("E21").Value = "  +3.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.735"
$ws.Range("E22").Value = "  +0.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.56"
$ws.Range("E23").Value = "  +2.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.26"
$ws.Range("E24").Value = "  +2.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  +6.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.61"
$ws.Range("E26").Value = "  +0.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.79"
$ws.Range("E28").Value = "  +5.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.42"
$ws.Range("E29").Value = "  +4.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.69"
$ws.Range("E30").Value = "  +0.61%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.19"
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("E33").Value = "  +6.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "27.39"
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("D35").Value = "0.0₃0840"
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.14"
$ws.Range("E37").Value = "  +2.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.32"
$ws.Range("E38").Value = "  -0.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.25"
$ws.Range("E39").Value = "  -3.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "9.33"
$ws.Range("E40").Value = "  +5.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "50.71"
$ws.Range("E41").Value = "  -1.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "450.36"
$ws.Range("E42").Value = "  +2.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.289"
$ws.Range("E43").Value = "  +1.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0368"
$ws.Range("E44").Value = "  -0.91%  "
$ws.Range("D45").Value = "2.842.36"
$ws.Range("E45").Value = "  -2.33%  "
$ws.Range("E46").Value = "  +1.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.51"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.82"
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.12"
$ws.Range("E50").Value = "  +4.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.26"
$ws.Range("E51").Value = "  +3.30%  "
